# NationalTreeSales.xlsx — "Fixed xlsx to be proper order"
#
# The sheet's H/I/J columns were saved out of order (Real trees / Fake
# trees / Avg home size landed as H/I/J instead of the intended
# Fake trees / Avg home size / Real trees). This reorders the three
# columns by rotating their contents left: new H = old I, new I = old J,
# new J = old H (for both the header row and all 17 data rows), and
# moves the active selection from I19 to L5, matching the saved sheet
# view of the fixed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1
$lastRow = 17

$colH = 8   # H
$colI = 9   # I
$colJ = 10  # J

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldH = $ws.Cells.Item($r, $colH).Value2
    $oldI = $ws.Cells.Item($r, $colI).Value2
    $oldJ = $ws.Cells.Item($r, $colJ).Value2

    $ws.Cells.Item($r, $colH).Value2 = $oldI
    $ws.Cells.Item($r, $colI).Value2 = $oldJ
    $ws.Cells.Item($r, $colJ).Value2 = $oldH
}

# Move the saved selection to match the re-saved sheet view.
[void]$ws.Range("L5").Select()
